# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" - row 2 (F2/G2), row 3 (F3), row 4 (F4)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 408
$ws1.Range("G2").Value = "已售罄"
$ws1.Range("F3").Value = 2330
$ws1.Range("F4").Value = 112

# Sheet "全部类型" - row 2 (F2/G2), row 7 (F7), row 8 (F8)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 408
$ws4.Range("G2").Value = "已售罄"
$ws4.Range("F7").Value = 2330
$ws4.Range("F8").Value = 112
